# Re-process the curated-dimensions metadata sheet.
# Rows 2-4 describe, per data column (A..G), the semantic annotations used
# during publication; row 5 (mapping-file pointers for the old
# "aragon"/"segunda-residencia" dimension mappings) is no longer needed and
# is dropped entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: iaest-measure / sdmx-dimension annotation
$ws.Range("A2").Value = "iaest-measure:numero-hogares"
$ws.Range("B2").Value = "sdmx-dimension:refArea"
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "null"
$ws.Range("F2").Value = "iaest-measure:segunda-residencia"
$ws.Range("G2").Value = "null"

# Row 3: measure / dimension kind
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "dim"
$ws.Range("C3").Value = "dim"
$ws.Range("D3").Value = "dim"
$ws.Range("E3").Value = "null"
$ws.Range("F3").Value = "medida"
$ws.Range("G3").Value = "null"

# Row 4: data type / URI template
$ws.Range("A4").Value = "xsd:int"
$ws.Range("B4").Value = "URI-Municipio"
$ws.Range("C4").Value = "URI-Provincia"
$ws.Range("D4").Value = "URI-Comunidad"
$ws.Range("E4").Value = "null"
$ws.Range("F4").Value = "xsd:int"
$ws.Range("G4").Value = "null"

# Row 5 (old mapping-file references) is obsolete now that the curated
# dimensions no longer rely on external mapping workbooks.
$ws.Range("A5:G5").EntireRow.Delete()
